# worklog.xlsx update
# - Remove the blank spacer row (old row 5) so the existing entries shift up
#   by one row (old row 6 -> new row 5, ... old row 11 -> new row 10).
# - Re-fit the row heights for the reflowed "Work Done" / "Bugs" text.
# - Fill in the (now un-blank) last existing row (new row 10) with the new
#   work-log entry text and a new "company name" column entry.
# - Append a brand new trailing row (row 11) for the next day.
# - Move the active selection to the new last cell, E11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the empty row 5 gap - everything below shifts up by one row.
$ws.Rows.Item(5).Delete()

# Row 5 (was row 6): height unchanged (60), content unchanged.

# Row 6 (was row 7): height 90 -> 75.
$ws.Rows.Item(6).RowHeight = 75

# Row 7 (was row 8): height 75 -> 45.
$ws.Rows.Item(7).RowHeight = 45

# Rows 8 and 9 (was 9 and 10): unchanged content/height (Holiday rows).

# Row 10 (was row 11, previously only had A/B filled in): add the new
# "Work Done" text (general/bottom aligned, wrapped) and the new
# company-name column entry (left/top aligned, wrapped, like the other
# entries in column E/C).
$ws.Cells.Item(10, 3).Value = "Fixed buggy scraping of pages where source code of the site does not load. Cleaned scraped text data (date,time and article). document and add comments to the code.`nkeep code in PEP-8 format."
$ws.Cells.Item(10, 3).HorizontalAlignment = 1
$ws.Cells.Item(10, 3).VerticalAlignment = -4107
$ws.Cells.Item(10, 3).WrapText = $true

$ws.Cells.Item(10, 5).Value = "Code works for only one company at a time"
$ws.Cells.Item(10, 5).HorizontalAlignment = -4131
$ws.Cells.Item(10, 5).VerticalAlignment = -4160
$ws.Cells.Item(10, 5).WrapText = $true

$ws.Rows.Item(10).RowHeight = 60

# Row 11 (brand new): next day's row with just the index and the date filled
# in, ready for future entries. Match column B's existing date formatting.
$ws.Cells.Item(11, 1).Value = 7
$ws.Cells.Item(11, 2).Value = 43641
$ws.Cells.Item(11, 2).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(11, 2).HorizontalAlignment = -4131
$ws.Cells.Item(11, 2).VerticalAlignment = -4160
$ws.Cells.Item(11, 2).WrapText = $true

# Update the selected cell to match where the user left off editing.
$ws.Range("E11").Select() | Out-Null
